$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# ---- Upper Section (rows 1-17) ----
$ws.Range("B2").Value = 0
$ws.Range("I3").Value = 15
$ws.Range("H9").Value = "ROUND THICKNESSES TO NEAREST 0.5CM!!"

$ws.Range("E13").Value = "Radiant Heat in (W)"
$ws.Range("F13").Value = 0

$ws.Range("J4").Value = "# not used for now"

$ws.Range("H10").Value = "DO NOT CHANGE LAYOUT"
$ws.Range("H10").Font.Size = 14
$ws.Range("H10").Font.Bold = $true
$ws.Range("H10").Font.Color = 255

$ws.Range("D11").Value = "Insulating plaster"
$ws.Range("B12").Value = 10

$ws.Range("A15").Value = "Radiant Heat in (W)"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = "Radiant Heat in (W)"
$ws.Range("D15").Value = 0

# ---- Middle Section (rows 19-35) ----
$ws.Range("B20").Value = 0
$ws.Range("D29").Value = "Insulating plaster"

$ws.Range("E31").Value = "Radiant Heat in (W)"
$ws.Range("F31").Value = 0

$ws.Range("A33").Value = "Radiant Heat in (W)"
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = "Radiant Heat in (W)"
$ws.Range("D33").Value = 0

# ---- Lower Section (rows 37-53) ----
$ws.Range("B38").Value = 0

$ws.Range("G49").Value = "Radiant Heat in (W)"
$ws.Range("H49").Value = 0

$ws.Range("A51").Value = "Radiant Heat in (W)"
$ws.Range("B51").Value = 0
$ws.Range("C51").Value = "Radiant Heat in (W)"
$ws.Range("D51").Value = 10000
$ws.Range("E51").Value = "Radiant Heat in (W)"
$ws.Range("F51").Value = 0

